# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps for the first data row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-21 05:05:40"
$wsZh.Range("H2").Value = "2016-03-21 05:06:19"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-21 05:05:49"
$wsDe.Range("H2").Value = "2016-03-21 05:06:33"
